# "New user added in excel" -- append a new "NewUser" worksheet (after the
# last existing sheet, "ValidationStrings") containing an Email/Password/EID
# header row plus one data row, then select it so it becomes the workbook's
# active sheet/tab.

$wb = $excel.ActiveWorkbook

# Add the new sheet after the current last sheet so it lands at the end of
# the tab strip (Worksheets.Add defaults to inserting before the active
# sheet, which is not what we want here).
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "NewUser"

# Header row (bold, matching the workbook's existing bold style).
$newSheet.Range("A1").Value = "Email"
$newSheet.Range("B1").Value = "Password"
$newSheet.Range("C1").Value = "EID"
$newSheet.Range("A1:C1").Font.Bold = $true

# Data row -- the user's first typed values, then corrected to the final
# email/EID (password stays the same), matching the edit history baked into
# the shared-string table.
$newSheet.Range("A2").Value = "TestPF1221+11082021150517@gmail.com"
$newSheet.Range("B2").Value = "pfqa_123"
$newSheet.Range("C2").Value = "TestPF1221_11082021150517"
$newSheet.Range("A2").Value = "TestPF1221+11082021152024@gmail.com"
$newSheet.Range("C2").Value = "TestPF1221_11082021152024"

$newSheet.PageSetup.Orientation = 1

# Leave the selection on B1 and make NewUser the active sheet/tab.
[void]$newSheet.Range("B1").Select()
$newSheet.Activate()
